# Daily attendance processing - normalize "Recorded By" (column G) entries
# so that a leading "System, " marker is moved to the end of the
# comma-separated list of recorders instead of the front.
#   e.g. "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#        "System, backup@backdoor.com, system" -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$prefix = "System, "

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $v = $cell.Text

    if ($v -ne $null -and $v.StartsWith($prefix)) {
        $rest = $v.Substring($prefix.Length)
        $cell.Value = $rest + ", System"
    }
}
